$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.825306666666667
$ws.Range("H2").Value = 20.47592
$ws.Range("I2").Value = 0.0939724583253512
$ws.Range("J2").Value = 0.09397245832535123
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.211539
$ws.Range("N2").Value = 0.634617
$ws.Range("O2").Value = 0.00536521120146958
$ws.Range("P2").Value = 0.005365211201469579
$ws.Range("Q2").Value = 1.44381854696
$ws.Range("R2").Value = 12.99436692264
$ws.Range("S2").Value = 0.0005041820860368075
$ws.Range("T2").Value = 0.0005041820860368076
$ws.Range("G3").Value = 6.825306666666667
$ws.Range("H3").Value = 20.47592
$ws.Range("I3").Value = 0.0939724583253512
$ws.Range("J3").Value = 0.09397245832535123
$ws.Range("O3").Value = 0.3044062438564017
$ws.Range("P3").Value = 0.3044062438564017
$ws.Range("Q3").Value = 81.91800176848888
$ws.Range("R3").Value = 737.2620159163999
$ws.Range("S3").Value = 0.0286058030647724
$ws.Range("T3").Value = 0.02860580306477241
$ws.Range("G4").Value = 6.825306666666667
$ws.Range("H4").Value = 20.47592
$ws.Range("I4").Value = 0.0939724583253512
$ws.Range("J4").Value = 0.09397245832535123
$ws.Range("M4").Value = 22.63137
$ws.Range("N4").Value = 67.89411
$ws.Range("O4").Value = 0.5739938253872932
$ws.Range("P4").Value = 0.573993825387293
$ws.Range("Q4").Value = 154.4660405368
$ws.Range("R4").Value = 1390.1943648312
$ws.Range("S4").Value = 0.05393961083521632
$ws.Range("T4").Value = 0.05393961083521633
$ws.Range("G5").Value = 6.825306666666667
$ws.Range("H5").Value = 20.47592
$ws.Range("I5").Value = 0.0939724583253512
$ws.Range("J5").Value = 0.09397245832535123
$ws.Range("M5").Value = 4.582890666666667
$ws.Range("N5").Value = 13.748672
$ws.Range("O5").Value = 0.1162347195548357
$ws.Range("P5").Value = 0.1162347195548357
$ws.Range("Q5").Value = 31.27963421980445
$ws.Range("R5").Value = 281.51670797824
$ws.Range("S5").Value = 0.01092286233932568
$ws.Range("T5").Value = 0.01092286233932568
$ws.Range("I6").Value = 0.5190671349373497
$ws.Range("J6").Value = 0.5190671349373498
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.211539
$ws.Range("N6").Value = 0.634617
$ws.Range("O6").Value = 0.00536521120146958
$ws.Range("P6").Value = 0.005365211201469579
$ws.Range("Q6").Value = 7.975089402739999
$ws.Range("R6").Value = 71.77580462466
$ws.Range("S6").Value = 0.00278490480668059
$ws.Range("T6").Value = 0.00278490480668059
$ws.Range("I7").Value = 0.5190671349373497
$ws.Range("J7").Value = 0.5190671349373498
$ws.Range("O7").Value = 0.3044062438564017
$ws.Range("P7").Value = 0.3044062438564017
$ws.Range("S7").Value = 0.1580072768555826
$ws.Range("T7").Value = 0.1580072768555827
$ws.Range("I8").Value = 0.5190671349373497
$ws.Range("J8").Value = 0.5190671349373498
$ws.Range("M8").Value = 22.63137
$ws.Range("N8").Value = 67.89411
$ws.Range("O8").Value = 0.5739938253872932
$ws.Range("P8").Value = 0.573993825387293
$ws.Range("Q8").Value = 853.2100419141999
$ws.Range("R8").Value = 7678.890377227799
$ws.Range("S8").Value = 0.2979413304155116
$ws.Range("T8").Value = 0.2979413304155116
$ws.Range("I9").Value = 0.5190671349373497
$ws.Range("J9").Value = 0.5190671349373498
$ws.Range("M9").Value = 4.582890666666667
$ws.Range("N9").Value = 13.748672
$ws.Range("O9").Value = 0.1162347195548357
$ws.Range("P9").Value = 0.1162347195548357
$ws.Range("Q9").Value = 172.7764752109511
$ws.Range("R9").Value = 1554.98827689856
$ws.Range("S9").Value = 0.0603336228595749
$ws.Range("T9").Value = 0.0603336228595749
$ws.Range("G10").Value = 15.52625766666667
$ws.Range("H10").Value = 46.578773
$ws.Range("I10").Value = 0.2137692374549467
$ws.Range("J10").Value = 0.2137692374549468
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.211539
$ws.Range("N10").Value = 0.634617
$ws.Range("O10").Value = 0.00536521120146958
$ws.Range("P10").Value = 0.005365211201469579
$ws.Range("Q10").Value = 3.284409020549
$ws.Range("R10").Value = 29.559681184941
$ws.Range("S10").Value = 0.001146917107322891
$ws.Range("T10").Value = 0.001146917107322891
$ws.Range("G11").Value = 15.52625766666667
$ws.Range("H11").Value = 46.578773
$ws.Range("I11").Value = 0.2137692374549467
$ws.Range("J11").Value = 0.2137692374549468
$ws.Range("O11").Value = 0.3044062438564017
$ws.Range("P11").Value = 0.3044062438564017
$ws.Range("Q11").Value = 186.3476712640038
$ws.Range("R11").Value = 1677.129041376035
$ws.Range("S11").Value = 0.06507269062570756
$ws.Range("T11").Value = 0.06507269062570757
$ws.Range("G12").Value = 15.52625766666667
$ws.Range("H12").Value = 46.578773
$ws.Range("I12").Value = 0.2137692374549467
$ws.Range("J12").Value = 0.2137692374549468
$ws.Range("M12").Value = 22.63137
$ws.Range("N12").Value = 67.89411
$ws.Range("O12").Value = 0.5739938253872932
$ws.Range("P12").Value = 0.573993825387293
$ws.Range("Q12").Value = 351.38048196967
$ws.Range("R12").Value = 3162.42433772703
$ws.Range("S12").Value = 0.1227022223568895
$ws.Range("T12").Value = 0.1227022223568895
$ws.Range("G13").Value = 15.52625766666667
$ws.Range("H13").Value = 46.578773
$ws.Range("I13").Value = 0.2137692374549467
$ws.Range("J13").Value = 0.2137692374549468
$ws.Range("M13").Value = 4.582890666666667
$ws.Range("N13").Value = 13.748672
$ws.Range("O13").Value = 0.1162347195548357
$ws.Range("P13").Value = 0.1162347195548357
$ws.Range("Q13").Value = 71.15514134882844
$ws.Range("R13").Value = 640.396272139456
$ws.Range("S13").Value = 0.02484740736502682
$ws.Range("T13").Value = 0.02484740736502682
$ws.Range("G14").Value = 12.579035
$ws.Range("H14").Value = 37.737105
$ws.Range("I14").Value = 0.1731911692823522
$ws.Range("J14").Value = 0.1731911692823523
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.211539
$ws.Range("N14").Value = 0.634617
$ws.Range("O14").Value = 0.00536521120146958
$ws.Range("P14").Value = 0.005365211201469579
$ws.Range("Q14").Value = 2.660956484865
$ws.Range("R14").Value = 23.948608363785
$ws.Range("S14").Value = 0.0009292072014292904
$ws.Range("T14").Value = 0.0009292072014292905
$ws.Range("G15").Value = 12.579035
$ws.Range("H15").Value = 37.737105
$ws.Range("I15").Value = 0.1731911692823522
$ws.Range("J15").Value = 0.1731911692823523
$ws.Range("O15").Value = 0.3044062438564017
$ws.Range("P15").Value = 0.3044062438564017
$ws.Range("Q15").Value = 150.9748150084416
$ws.Range("R15").Value = 1358.773335075975
$ws.Range("S15").Value = 0.05272047331033906
$ws.Range("T15").Value = 0.05272047331033908
$ws.Range("G16").Value = 12.579035
$ws.Range("H16").Value = 37.737105
$ws.Range("I16").Value = 0.1731911692823522
$ws.Range("J16").Value = 0.1731911692823523
$ws.Range("M16").Value = 22.63137
$ws.Range("N16").Value = 67.89411
$ws.Range("O16").Value = 0.5739938253872932
$ws.Range("P16").Value = 0.573993825387293
$ws.Range("Q16").Value = 284.68079532795
$ws.Range("R16").Value = 2562.12715795155
$ws.Range("S16").Value = 0.09941066177967561
$ws.Range("T16").Value = 0.09941066177967563
$ws.Range("G17").Value = 12.579035
$ws.Range("H17").Value = 37.737105
$ws.Range("I17").Value = 0.1731911692823522
$ws.Range("J17").Value = 0.1731911692823523
$ws.Range("M17").Value = 4.582890666666667
$ws.Range("N17").Value = 13.748672
$ws.Range("O17").Value = 0.1162347195548357
$ws.Range("P17").Value = 0.1162347195548357
$ws.Range("Q17").Value = 57.64834209717333
$ws.Range("R17").Value = 518.8350788745599
$ws.Range("S17").Value = 0.02013082699090829
$ws.Range("T17").Value = 0.02013082699090829

Write-Host "Updated 182 cells"
